$wb = $excel.ActiveWorkbook

# Add a new worksheet after the existing "Statistics" sheet
$statsSheet = $wb.Worksheets.Item("Statistics")
$newSheet = $wb.Worksheets.Add($null, $statsSheet)
$newSheet.Name = "MinVarPortfolio"

# Header row: Asset, Weight (set Weight first so it is registered earlier in the shared-strings table)
$newSheet.Range("B1").Value = "Weight"
$newSheet.Range("A1").Value = "Asset"

# Data rows
$assets = @("NASDAQ", "S&P500", "DAX", "WIG20")
for ($i = 0; $i -lt $assets.Length; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 1).Value = $assets[$i]
    $newSheet.Cells.Item($row, 2).Value = 0.25
}

# Apply the bold/centered/bordered "label" formatting (used for header-like cells on
# the Statistics sheet) to the header row and the asset-name column, reusing the
# existing style via a format-only copy/paste so no new style entries are created.
$statsSheet.Range("A1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)
$newSheet.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the original sheet as the active one (Worksheets.Add activates the new sheet)
$statsSheet.Activate()
